$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109, shifting existing rows 109-170 down to 110-171
$ws.Rows.Item(109).Insert()

# Populate the new row 109 with the weekly data point
$ws.Cells.Item(109, 1).Value = 7
$ws.Cells.Item(109, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(109, 3).Value = "Ñuble"
$ws.Cells.Item(109, 4).Value = 44488
$ws.Cells.Item(109, 5).Value = 16
$ws.Cells.Item(109, 6).Value = 100112009
$ws.Cells.Item(109, 7).Value = "Acelga"
$ws.Cells.Item(109, 8).Value = "Sin especificar"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 120
$ws.Cells.Item(109, 11).Value = 350
$ws.Cells.Item(109, 12).Value = 400
$ws.Cells.Item(109, 13).Value = 375
$ws.Cells.Item(109, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(109, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(109, 16).Value = 375
$ws.Cells.Item(109, 17).Value = 1
$ws.Cells.Item(109, 18).Value = "Hortaliza"
